$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1379
$ws.Range("I2").Value = 3638
$ws.Range("J2").Value = 14741
$ws.Range("K2").Value = 76
$ws.Range("L2").Value = 4069
$ws.Range("M2").Value = 236
$ws.Range("N2").Value = 2624
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 52
$ws.Range("Q2").Value = 28
$ws.Range("R2").Value = 183
$ws.Range("S2").Value = 1585
$ws.Range("T2").Value = 2547
$ws.Range("U2").Value = 199
$ws.Range("V2").Value = 23409
$ws.Range("W2").Value = 14
$ws.Range("X2").Value = 23024
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 358
$ws.Range("AA2").Value = 146
